$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.975.20'
$ws.Range("E2").Value = '  -3.56%  '
$ws.Range("D3").Value = '3.504.93'
$ws.Range("E3").Value = '  -2.81%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '586.53'
$ws.Range("E5").Value = '  -3.34%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '132.31'
$ws.Range("E6").Value = '  -5.76%  '
$ws.Range("D7").Value = '3.505.09'
$ws.Range("E7").Value = '  -2.82%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("E9").Value = '  -1.44%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.125'
$ws.Range("E10").Value = '  -2.38%  '
$ws.Range("E11").Value = '  -1.80%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.386'
$ws.Range("E12").Value = '  -2.69%  '
$ws.Range("D13").Value = '4.096.08'
$ws.Range("E13").Value = '  -3.04%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.90'
$ws.Range("E14").Value = '  -3.00%  '
$ws.Range("E15").Value = '  -4.51%  '
$ws.Range("E16").Value = '  +0.55%  '
$ws.Range("D17").Value = '3.505.77'
$ws.Range("E17").Value = '  -2.83%  '
$ws.Range("D18").Value = '64.059.69'
$ws.Range("E18").Value = '  -3.47%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.01'
$ws.Range("E19").Value = '  -2.17%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.49'
$ws.Range("E20").Value = '  -1.87%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.68'
$ws.Range("E21").Value = '  -4.40%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '391.69'
$ws.Range("E22").Value = '  -2.03%  '
$ws.Range("E23").Value = '  -2.78%  '
$ws.Range("D24").Value = '3.644.59'
$ws.Range("E24").Value = '  -2.87%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '72.93'
$ws.Range("E25").Value = '  -3.05%  '
$ws.Range("E26").Value = '  +0.23%  '
$ws.Range("E27").Value = '  -6.74%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.60'
$ws.Range("E28").Value = '  -3.69%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.49'
$ws.Range("E29").Value = '  -8.82%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.999'
$ws.Range("E30").Value = '  +0.01%  '
$ws.Range("E31").Value = '  -4.10%  '
$ws.Range("E32").Value = '  -5.45%  '
$ws.Range("D33").Value = '3.509.91'
$ws.Range("E33").Value = '  -2.85%  '
$ws.Range("E34").Value = '  +0.04%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '23.90'
$ws.Range("E35").Value = '  -3.51%  '
$ws.Range("E36").Value = '  -3.68%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.32'
$ws.Range("E37").Value = '  -2.25%  '
$ws.Range("E38").Value = '  -5.20%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.96'
$ws.Range("E39").Value = '  -2.33%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '166.98'
$ws.Range("E40").Value = '  -1.15%  '
$ws.Range("E41").Value = '  -4.40%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '27.39'
$ws.Range("E42").Value = '  +2.25%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.814'
$ws.Range("E43").Value = '  -3.97%  '
$ws.Range("E44").Value = '  -0.18%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '41.94'
$ws.Range("E45").Value = '  -2.92%  '
$ws.Range("E46").Value = '  -6.19%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.40'
$ws.Range("E47").Value = '  -4.48%  '
$ws.Range("E48").Value = '  -5.40%  '
$ws.Range("D49").Value = '2.449.61'
$ws.Range("E49").Value = '  -0.83%  '
$ws.Range("E50").Value = '  -2.59%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.901'
$ws.Range("E51").Value = '  -1.70%  '
